$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1503136666666667
$ws.Range("N2").Value = 0.450941
$ws.Range("O2").Value = 0.008029526741163598
$ws.Range("P2").Value = 0.008029526741163598
$ws.Range("Q2").Value = 1.568409374325556
$ws.Range("R2").Value = 14.11568436893
$ws.Range("S2").Value = 0.007797547272960633
$ws.Range("T2").Value = 0.007797547272960634
$ws.Range("O3").Value = 0.8389317081486641
$ws.Range("P3").Value = 0.8389317081486641
$ws.Range("S3").Value = 0.8146942981756413
$ws.Range("T3").Value = 0.8146942981756414
$ws.Range("M4").Value = 2.758466666666667
$ws.Range("N4").Value = 8.2754
$ws.Range("O4").Value = 0.1473530807662759
$ws.Range("P4").Value = 0.1473530807662759
$ws.Range("Q4").Value = 28.78251242688889
$ws.Range("R4").Value = 259.0426118419999
$ws.Range("S4").Value = 0.143095932067961
$ws.Range("T4").Value = 0.143095932067961
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1064366666666667
$ws.Range("N5").Value = 0.31931
$ws.Range("O5").Value = 0.005685684343896315
$ws.Range("P5").Value = 0.005685684343896314
$ws.Range("Q5").Value = 1.110586079588889
$ws.Range("R5").Value = 9.995274716299999
$ws.Range("S5").Value = 0.005521420362595239
$ws.Range("T5").Value = 0.005521420362595239
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1503136666666667
$ws.Range("N6").Value = 0.450941
$ws.Range("O6").Value = 0.008029526741163598
$ws.Range("P6").Value = 0.008029526741163598
$ws.Range("Q6").Value = 0.046660669034
$ws.Range("R6").Value = 0.419946021306
$ws.Range("S6").Value = 0.0002319794682029637
$ws.Range("T6").Value = 0.0002319794682029637
$ws.Range("O7").Value = 0.8389317081486641
$ws.Range("P7").Value = 0.8389317081486641
$ws.Range("S7").Value = 0.0242374099730227
$ws.Range("T7").Value = 0.0242374099730227
$ws.Range("M8").Value = 2.758466666666667
$ws.Range("N8").Value = 8.2754
$ws.Range("O8").Value = 0.1473530807662759
$ws.Range("P8").Value = 0.1473530807662759
$ws.Range("Q8").Value = 0.8562887396
$ws.Range("R8").Value = 7.706598656399999
$ws.Range("S8").Value = 0.004257148698314868
$ws.Range("T8").Value = 0.004257148698314869
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1064366666666667
$ws.Range("N9").Value = 0.31931
$ws.Range("O9").Value = 0.005685684343896315
$ws.Range("P9").Value = 0.005685684343896314
$ws.Range("Q9").Value = 0.03304028294
$ws.Range("R9").Value = 0.29736254646
$ws.Range("S9").Value = 0.0001642639813010756
$ws.Range("T9").Value = 0.0001642639813010756
